$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: fix E35 (to numeric phone number) and F35 (timestamp) ---
$ws.Range("E35").Value = 917990747606
$ws.Range("F35").Value = 45964.77464993056
$ws.Range("F35").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 36: [Call Started] ---
$ws.Range("A36").Value = '[Call Started]'
$ws.Range("D36").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E36").Value = 917990747606
$ws.Range("F36").Value = 45964.78596905093
$ws.Range("F36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(36).EntireRow.AutoFit()

# --- Row 37: [Intro response] ---
$ws.Range("A37").Value = '[Intro response]'
$ws.Range("B37").Value = 'Yes.'
$ws.Range("C37").Value = 'neutral'
$text37 = @'
Here are our latest offers:
- Laptop Pro 
- Smart watch 
- Bluetooth Earbuds
Which product would you like to purchase?
'@
$ws.Range("D37").Value = $text37
$ws.Range("E37").Value = 917990747606
$ws.Range("F37").Value = 45964.78619211805
$ws.Range("F37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(37).EntireRow.AutoFit()

# --- Row 38: [Product match] ---
$ws.Range("A38").Value = '[Product match]'
$ws.Range("B38").Value = 'Bluetooth earbuds.'
$ws.Range("C38").Value = 'neutral'
$ws.Range("D38").Value = 'Great choice! I’ve sent the link of Bluetooth Earbuds to your phone number ending with 7606. Thank you for your time! I really appreciate it.'
$ws.Range("E38").Value = 917990747606
$ws.Range("F38").Value = 45964.78638554398
$ws.Range("F38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(38).EntireRow.AutoFit()

# --- Row 39: [Call Started] ---
$ws.Range("A39").Value = '[Call Started]'
$ws.Range("D39").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E39").Value = 917990747606
$ws.Range("F39").Value = 45964.82465596065
$ws.Range("F39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(39).EntireRow.AutoFit()

# --- Row 40: [Persuasion check] ---
$ws.Range("A40").Value = '[Persuasion check]'
$ws.Range("B40").Value = 'No.'
$ws.Range("C40").Value = 'neutral'
$ws.Range("D40").Value = 'I completely understand! But before you go — we’re giving a 20% discount just for today. Would you like to take a quick look?'
$ws.Range("E40").Value = 917990747606
$ws.Range("F40").Value = 45964.82489373843
$ws.Range("F40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(40).EntireRow.AutoFit()

# --- Row 41: [Persuasion check] ---
$ws.Range("A41").Value = '[Persuasion check]'
$ws.Range("B41").Value = 'Uh, just not interested.'
$ws.Range("C41").Value = 'neutral'
$ws.Range("D41").Value = 'Totally fair! But if I may — we’re offering free delivery on all products this week. Can I share a few top deals?'
$ws.Range("E41").Value = 917990747606
$ws.Range("F41").Value = 45964.82510863426
$ws.Range("F41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(41).EntireRow.AutoFit()

# --- Row 42: [Call Started] (E42 keeps the leading '+' as text) ---
$ws.Range("A42").Value = '[Call Started]'
$ws.Range("D42").Value = 'Hi there! I’m your sales agent from Creer Infotech. I’ve reached out to share some exciting offers on our latest products. Can I take a few minutes to tell you about them?'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '+917990747606'
$ws.Range("E42").Style = "Normal"
$ws.Range("F42").Value = 45964.82681271759
$ws.Range("F42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Rows(42).EntireRow.AutoFit()
